# "added foreign key + check"
# Adds a new "Team" column (J) to Sheet1: header "Team" in J1 and the
# literal value "nor" for every data row (J2:J89), matching the header
# styling already used for B1:I1 (bold, centered/top, boxed) but with a
# left+right border only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header + fill down the constant "nor" for all data rows.
$ws.Range("J1").Value = "Team"
$ws.Range("J2:J89").Value = "nor"

# Style the new header cell like the others (bold, centered, top-aligned)
# but boxed only on the left/right edges.
$hdr = $ws.Range("J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.Item(7).LineStyle = 1   # xlEdgeLeft  -> xlContinuous
$hdr.Borders.Item(10).LineStyle = 1  # xlEdgeRight -> xlContinuous

# Reflect the new data extent in the active selection.
$null = $ws.Range("A2:J89").Select()
